$wb = $excel.ActiveWorkbook

# --- Switch to the "Repayment schedule" sheet (it becomes the active tab) ---
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate() | Out-Null

# --- Insert a new blank column before column N (shifts N:P -> O:Q) ---
$ws.Columns("N").Insert() | Out-Null

# Give the freshly inserted column the same width as its neighbour to its
# left (column M), matching what Excel does visually when a column is
# inserted in the middle of a formatted table.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# --- Leave the selection on S5, as last left by the editor ---
$ws.Range("S5").Select() | Out-Null
